$wb = $excel.ActiveWorkbook

# Rename the two test-case sheets: tc_005 -> tc_006, tc_016 -> tc_007
$wb.Worksheets.Item("tc_005").Name = "tc_006"
$wb.Worksheets.Item("tc_016").Name = "tc_007"

$wsTc006 = $wb.Worksheets.Item("tc_006")
$wsTc007 = $wb.Worksheets.Item("tc_007")

# tc_007 (previously tc_016) keeps its selection at F10 but is no longer the
# tab-selected/active sheet.
$wsTc007.Select() | Out-Null
$wsTc007.Range("F10").Select() | Out-Null

# tc_006 (previously tc_005) becomes the active/tab-selected sheet, with its
# selection moved from B2 to D10.
$wsTc006.Select() | Out-Null
$wsTc006.Range("D10").Select() | Out-Null
